# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    $rng = $ws.Range($Addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '63.070.82'
$ws.Range('E2').Value = '  -1.52%  '
Set-TextValue 'D3' '3.056.92'
$ws.Range('E3').Value = '  -3.00%  '
$ws.Range('E4').Value = '  -0.29%  '
Set-TextValue 'D5' '589.54'
$ws.Range('E5').Value = '  -0.45%  '
Set-TextValue 'D6' '152.52'
$ws.Range('E6').Value = '  +4.38%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  +3.07%  '
Set-TextValue 'D9' '3.058.78'
$ws.Range('E9').Value = '  -2.66%  '
$ws.Range('E10').Value = '  -3.60%  '
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('E12').Value = '  -0.05%  '
Set-TextValue 'D13' '0.0000240'
$ws.Range('E13').Value = '  -2.85%  '
Set-TextValue 'D14' '37.02'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('E15').Value = '  -1.89%  '
Set-TextValue 'D16' '3.564.16'
$ws.Range('E16').Value = '  -3.00%  '
Set-TextValue 'D17' '7.19'
$ws.Range('E17').Value = '  -1.23%  '
Set-TextValue 'D18' '63.140.17'
Set-TextValue 'D19' '3.058.54'
$ws.Range('E19').Value = '  -2.86%  '
Set-TextValue 'D20' '473.30'
$ws.Range('E20').Value = '  +1.20%  '
Set-TextValue 'D21' '14.61'
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('E22').Value = '  -2.33%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('E24').Value = '  +1.56%  '
Set-TextValue 'D25' '12.94'
$ws.Range('E25').Value = '  -0.46%  '
Set-TextValue 'D26' '81.10'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D29' '2.67'
$ws.Range('E29').Value = '  -1.43%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D30' '1.00'
$ws.Range('E30').Value = '  -0.11%  '
Set-TextValue 'D31' '7.26'
$ws.Range('E32').Value = '  -2.16%  '
Set-TextValue 'D33' '0.114'
$ws.Range('E33').Value = '  +2.75%  '
$ws.Range('E34').Value = '  -1.90%  '
Set-TextValue 'D35' '0.0₃0838'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  -2.23%  '
$ws.Range('E37').Value = '  -1.18%  '
Set-TextValue 'D38' '3.32'
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('E39').Value = '  -4.95%  '
Set-TextValue 'D40' '9.26'
$ws.Range('E40').Value = '  +0.58%  '
Set-TextValue 'D41' '50.35'
$ws.Range('E41').Value = '  -2.04%  '
Set-TextValue 'D42' '442.72'
$ws.Range('E42').Value = '  -4.48%  '
$ws.Range('E43').Value = '  -3.38%  '
Set-TextValue 'D44' '40.15'
$ws.Range('E44').Value = '  -0.08%  '
Set-TextValue 'D45' '0.0361'
$ws.Range('E45').Value = '  -2.89%  '
$ws.Range('E46').Value = '  +2.05%  '
Set-TextValue 'D47' '2.793.06'
$ws.Range('E47').Value = '  -4.59%  '
Set-TextValue 'D48' '131.10'
$ws.Range('E48').Value = '  +1.58%  '
$ws.Range('E49').Value = '  +0.05%  '
Set-TextValue 'D50' '25.08'
$ws.Range('E50').Value = '  +3.48%  '
$ws.Range('E51').Value = '  +0.04%  '
